# Append a freshly-scraped job listing as the new row 4, pushing the
# previously-fetched rows (old rows 4-12) down to rows 5-13, and refresh
# every "取得日時" (fetched-at) timestamp to the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-08 01:18:26"

# --- 1. shift existing data rows 4..12 down to 5..13 (bottom-up to avoid clobbering) ---
for ($r = 12; $r -ge 4; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 8; $c++) {
        $val = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($dest, $c).Value = $val
    }
}

# --- 2. write the new row 4 (freshly scraped listing) ---
$ws.Cells.Item(4, 1).Value = $newTimestamp
$ws.Cells.Item(4, 2).Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5314730"
$ws.Cells.Item(4, 7).Value = 178
$ws.Cells.Item(4, 8).Value = "★bot ◆ツール"

# --- 3. refresh the "取得日時" timestamp on every data row (2..13) ---
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 4. rebuild hyperlinks on column F (URL) so each points at the row's own URL.
# (Hyperlinks collection on a sub-range deletes sheet-wide in this host, so clear
# once up front and recreate every row's link to land in the right order/target.)
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 13; $r++) {
    $url = $ws.Cells.Item($r, 6).Value()
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url)
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
